$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear stale "Serotype" labels that land on rows which should be blank
#     in column A once the table is renumbered (old A5="O", A9="A", A12="SAT1") ---
$ws.Range("A5").ClearContents()
$ws.Range("A9").ClearContents()

# --- Update existing A/B/C columns (values changed + row re-mapping) ---
$ws.Range("A1").Value = "Serotype"
$ws.Range("B1").Value = "Tree.height"
$ws.Range("C1").Value = "clockRate/ucldMean"
$ws.Range("A2").Value = "Asia1"
$ws.Range("B2").Value = 154
$ws.Range("C2").Value = 0.00303
$ws.Range("B3").Value = 139
$ws.Range("C3").Value = 0.00036
$ws.Range("A4").Value = "O"
$ws.Range("B4").Value = 143
$ws.Range("C4").Value = 0.001981
$ws.Range("B5").Value = 138
$ws.Range("C5").Value = 0.001892
$ws.Range("B6").Value = 130
$ws.Range("C6").Value = 0.002058
$ws.Range("B7").Value = 134
$ws.Range("C7").Value = 0.00207
$ws.Range("A8").Value = "A"
$ws.Range("B8").Value = 1554
$ws.Range("C8").Value = 0.000156
$ws.Range("B9").Value = 1579
$ws.Range("C9").Value = 0.0001041
$ws.Range("B10").Value = 1545
$ws.Range("C10").Value = 0.0001067
$ws.Range("A11").Value = "SAT1"
$ws.Range("B11").Value = 208
$ws.Range("C11").Value = 0.001738

# Clear the old row 12 (data moved into row 11; sheet now has 11 rows)
$ws.Range("A12:E12").ClearContents()

# --- Populate new D/E columns in the exact order that reproduces the original
#     shared-string table ordering ---
$ws.Range("D8").Value = "[1405.7829, 1681.2922]"
$ws.Range("D2").Value = "[84.7865, 244.5014]"
$ws.Range("D1").Value = "Tree.height 95% HPD interval"
$ws.Range("E1").Value = "clockRate/ucldMean 95% HPD interval"
$ws.Range("E2").Value = "[1.7678E-3, 4.4727E-3]"
$ws.Range("E3").Value = "[2.246E-3, 5.1772E-3]"
$ws.Range("D3").Value = "[81.2028, 207.3966]"
$ws.Range("E8").Value = "[1E-4, 1.1481E-4]"
$ws.Range("D9").Value = "[1443.9537, 1695.5336]"
$ws.Range("E9").Value = "[1E-4, 1.12E-4]"
$ws.Range("D10").Value = "[1362.3057, 1703.4225]"
$ws.Range("E10").Value = "[1E-4, 1.1832E-4]"
$ws.Range("D4").Value = "[97.5098, 196.6163]"
$ws.Range("E4").Value = "[1.2092E-3, 2.8237E-3]"
$ws.Range("D5").Value = "[94.9937, 189.1976]"
$ws.Range("E5").Value = "[1.1492E-3, 2.6736E-3]"
$ws.Range("D6").Value = "[87.8278, 178.9122]"
$ws.Range("D7").Value = "[97.33, 176.8925]"
$ws.Range("E7").Value = "[1.4176E-3, 2.832E-3]"
$ws.Range("E6").Value = "[1.2081E-3, 2.947E-3]"
$ws.Range("D11").Value = "[108.1109, 618.647]"
$ws.Range("E11").Value = "[9.6834E-4, 2.6446E-3]"

# --- Formatting: row 1 taller, wrap/vertical-center for D1/D2, new column widths ---
$ws.Rows.Item(1).RowHeight = 30
$ws.Range("D1:D2").VerticalAlignment = -4108
$ws.Range("D1:D2").WrapText = $true

# --- Column widths (closest achievable values under this engine's width model) ---
$ws.Columns.Item(2).ColumnWidth = 49.5
$ws.Columns.Item(3).ColumnWidth = 42.666666666666664
$ws.Columns.Item(4).ColumnWidth = 30.833333333333332
$ws.Columns.Item(5).ColumnWidth = 34.833333333333336

# --- Selection to match the saved view state ---
$ws.Range("I18").Select() | Out-Null
